# Swap the deck's theme palette from the "Integral" colour set to the
# stock "Office Theme" colour set (the slide master / design's theme
# colours are changed in place; dk1/lt1 are already 000000/FFFFFF in
# both palettes so only the remaining ten roles actually move).
#
# ThemeColorScheme item order (1-based):
#   1 dk1   2 lt1   3 dk2   4 lt2   5 accent1  6 accent2
#   7 accent3  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
#
# RGB values below are packed as R + G*256 + B*65536 (the usual
# VBA/COM colour-value encoding), matching the target hex:
#   dk2=44546A lt2=E7E6E6 accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5
#   accent4=FFC000 accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

$p   = $ppt.ActivePresentation
$sm  = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000 (unchanged)
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF (unchanged)
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
